$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the existing row 145,
# pushing all subsequent rows (145-169) down by one (to 146-170).
$ws.Rows.Item(145).EntireRow.Insert()

# Populate the newly inserted row 145 with the new record's data.
$ws.Range("A145").Value = 4
$ws.Range("B145").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C145").Value = "Los Lagos"
$ws.Range("D145").Value = 44504
$ws.Range("E145").Value = 10
$ws.Range("F145").Value = 100112043
$ws.Range("G145").Value = "Pepino ensalada"
$ws.Range("H145").Value = "Sin especificar"
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 150
$ws.Range("K145").Value = 12000
$ws.Range("L145").Value = 12000
$ws.Range("M145").Value = 12000
$ws.Range("N145").Value = "$/caja 60 unidades"
$ws.Range("O145").Value = "Región de Arica y Parinacota"
$ws.Range("P145").Value = 200
$ws.Range("Q145").Value = 60
$ws.Range("R145").Value = "Hortaliza"
